# Update "Name of Algo" result values in column C (KNN imputed values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = -12.462
$ws.Range("C6").Value  = -12.445
$ws.Range("C7").Value  = -12.673
$ws.Range("C16").Value = -12.126
$ws.Range("C20").Value = -13.041
